$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before row 9 (current row 9 = "SPASMOFEN ..."),
# pushing it and everything below it down by one row.
$ws.Rows.Item(9).Insert()

# Copy the formatting (styles, borders, merges) of the row that is now
# row 10 (the old row 9 / SPASMOFEN row) into the newly-inserted blank row 9,
# so the new row looks exactly like the other item rows in the table.
$ws.Rows.Item(10).Copy()
$ws.Rows.Item(9).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# Fill the new row (item #3 - MAVILOR) with its data.
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "MAVILOR 5MG 30 TAB."
$ws.Range("H9").Value = "0:0"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "105.00"
$ws.Range("P9").Value = "105.0000"
$ws.Range("Q9").Value = "1:0"

# Renumber the item index (column A) of the rows that shifted down.
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9

# Update the grand-total cell (now on row 16) to include the new item.
$ws.Range("P16").Value = 398.76499999999999

# Update the generated timestamp in the footer (now on row 17).
$ws.Range("A17").Value = "Friday, 20 June, 2025 4:30 PM"
